$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hungary NB I")

# ---------------------------------------------------------------------------
# Rows 102-104: odds/results data re-shuffled between the three rows
# (the A-column running id stays tied to the row, everything else moves).
# ---------------------------------------------------------------------------

# Row 102 (becomes old row 103's data)
$ws.Range("B102").Value = 5470380
$ws.Range("F102").Value = "Puskas Academy"
$ws.Range("G102").Value = "Budapest Honved"
$ws.Range("H102").Value = 2
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = "H"
$ws.Range("K102").Value = 2
$ws.Range("L102").Value = 3.6
$ws.Range("M102").Value = 3.6
$ws.Range("N102").Value = 1.75
$ws.Range("O102").Value = 3.8
$ws.Range("P102").Value = 4.5
$ws.Range("Q102").Value = -0.75
$ws.Range("R102").Value = 2
$ws.Range("S102").Value = 1.85
$ws.Range("U102").Value = 2
$ws.Range("V102").Value = 1.85
$ws.Range("W102").Value = 0.75
$ws.Range("X102").Value = -1
$ws.Range("Z102").Value = 0.5
$ws.Range("AA102").Value = -0.5
$ws.Range("AB102").Value = 0.5
$ws.Range("AC102").Value = -0.5

# Row 103 (becomes old row 104's data)
$ws.Range("B103").Value = 5461530
$ws.Range("F103").Value = "Debreceni VSC"
$ws.Range("G103").Value = "Ujpest"
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 1.727
$ws.Range("L103").Value = 3.5
$ws.Range("M103").Value = 4
$ws.Range("N103").Value = 1.7
$ws.Range("O103").Value = 4
$ws.Range("P103").Value = 4.75
$ws.Range("R103").Value = 1.875
$ws.Range("S103").Value = 1.975
$ws.Range("U103").Value = 1.975
$ws.Range("V103").Value = 1.875
$ws.Range("W103").Value = 0.7
$ws.Range("Z103").Value = 0.875
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = -1
$ws.Range("AC103").Value = 0.875

# Row 104 (becomes old row 102's data)
$ws.Range("B104").Value = 5461531
$ws.Range("F104").Value = "Vasas SC"
$ws.Range("G104").Value = "MOL Fehervar FC"
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = "D"
$ws.Range("K104").Value = 6
$ws.Range("L104").Value = 4.5
$ws.Range("M104").Value = 1.4
$ws.Range("N104").Value = 5.75
$ws.Range("O104").Value = 4.75
$ws.Range("P104").Value = 1.5
$ws.Range("Q104").Value = 1
$ws.Range("R104").Value = 2.05
$ws.Range("S104").Value = 1.8
$ws.Range("U104").Value = 1.825
$ws.Range("V104").Value = 2.025
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = 3.75
$ws.Range("Z104").Value = 1.05
$ws.Range("AC104").Value = 1.025

# ---------------------------------------------------------------------------
# Rows 226-228: refreshed odds for upcoming (not yet played) fixtures.
# ---------------------------------------------------------------------------

# Row 226
$ws.Range("B226").Value = 6818320
$ws.Range("E226").Value = 45339.47916666666
$ws.Range("F226").Value = "Mezokovesd Zsory"
$ws.Range("G226").Value = "Kisvarda FC"
$ws.Range("K226").Value = 2.15
$ws.Range("L226").Value = 3.3
$ws.Range("M226").Value = 3.3
$ws.Range("N226").Value = 2.75
$ws.Range("O226").Value = 3.1
$ws.Range("P226").Value = 2.625
$ws.Range("Q226").Value = 0
$ws.Range("T226").Value = 2.25
$ws.Range("U226").Value = 2
$ws.Range("V226").Value = 1.85

# Row 227
$ws.Range("B227").Value = 6818318
$ws.Range("E227").Value = 45340.45833333334
$ws.Range("F227").Value = "MOL Fehervar FC"
$ws.Range("G227").Value = "Debreceni VSC"
$ws.Range("K227").Value = 2.6
$ws.Range("L227").Value = 3.4
$ws.Range("M227").Value = 2.55
$ws.Range("N227").Value = 2.6
$ws.Range("O227").Value = 3.4
$ws.Range("R227").Value = 1.975
$ws.Range("S227").Value = 1.875
$ws.Range("T227").Value = 2.75
$ws.Range("U227").Value = 2.025
$ws.Range("V227").Value = 1.825

# Row 228
$ws.Range("B228").Value = 6818319
$ws.Range("E228").Value = 45340.5625
$ws.Range("F228").Value = "Puskas Academy"
$ws.Range("G228").Value = "MTK Budapest"
$ws.Range("K228").Value = 1.65
$ws.Range("L228").Value = 3.8
$ws.Range("M228").Value = 5
$ws.Range("N228").Value = 1.533
$ws.Range("O228").Value = 4
$ws.Range("P228").Value = 6
$ws.Range("Q228").Value = -1
$ws.Range("R228").Value = 1.925
$ws.Range("S228").Value = 1.925
$ws.Range("U228").Value = 2
$ws.Range("V228").Value = 1.85

# ---------------------------------------------------------------------------
# Row 229 no longer exists in the updated dataset - delete it entirely.
# ---------------------------------------------------------------------------
$ws.Rows.Item(229).Delete()
